# Remove the "CSP01" row (row 2) from Sheet1.
# This deletes the entire row and shifts all rows below it up by one,
# which matches the target diff (dimension shrinks from A1:C28 to A1:C27).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(2).Delete()
